# Finalise Fatigue Module edit:
# - Recreate relations that had been removed during modularisation
# - Set 'FatigueUseCase' as rootclass for all concepts in the Fatigue Module
#
# This reproduces the target diff against xl/worksheets/sheet1.xml and
# xl/sharedStrings.xml using plain Excel object-model operations (row
# insertion + cell value / font assignment) instead of touching the
# underlying XML directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: make room ---------------------------------------------------
# Insert one row right after the current row 46 ("What's more..." bullet
# list). This pushes the old row 47 ("...add Time...") down to row 48 and
# the old row 48 ("...rearrange labels...") down to row 49.
$ws.Rows.Item(47).Insert()

# Insert 15 more rows starting at row 50 to make space for the new bullet
# list (rows 50-64), so that the existing bottom block (the three shaded
# rows that used to be E58:F60) lands exactly on rows 74-76.
$ws.Rows("50:64").Insert()

# --- Step 2: fill in the new cells ---------------------------------------
# (values are entered in the same order they were originally typed, so
# that the shared-strings table comes out in the same order)
$ws.Range("A66").Value = 'In Fatigue module:'
$ws.Range("A67").Value = 'Revise Failure concepts:'
$ws.Range("B67").Value = 'FailureComment, FailureOrigin, FailureReason (subclasses Featureless, Inclusion (subclass InclusionElements), Pore, FailureType (subclasses: SurfaceFailure, VolumeFailure)'
$ws.Range("A68").Value = 'Fractography'
$ws.Range("B68").Value = 'Currently has superclass "Measurement", but isnt it a method, and a FractographyMeasurement is the measurement itself?'
$ws.Range("A50").Value = '…revise all concepts for suitability as subclasses of:'
$ws.Range("A51").Value = 'Device [from oie/manufacturing]'
$ws.Range("A52").Value = 'Characterisation Machine [chameo]'
$ws.Range("A53").Value = 'Detector [chameo]'
$ws.Range("A54").Value = '"Calibration Data" [chameo]'
$ws.Range("A56").Value = '"Primary Data" [chameo]'
$ws.Range("A57").Value = '"Raw Data" [chameo]'
$ws.Range("A58").Value = '"Secondary Data" [chameo]'
$ws.Range("A59").Value = '"Characterisation Property" [chameo]'
$ws.Range("A47").Value = '…"Device" is double, 1x from oie/manufacturing and 1x from emmo beta4, the classes have different meanings, IRIs, and subclasses, but the same label.'
$ws.Range("B69").Value = 'same as Fractography'
$ws.Range("A69").Value = 'FatigueTesting and subclasses'
$ws.Range("A61").Value = '"Characterisation Method" [chameo]'
$ws.Range("A60").Value = '"Measurement Process" [chameo]'
$ws.Range("A55").Value = '"Calibration Measurement" [chameo]'
$ws.Range("A62").Value = '"Characterisation Workflow" [chameo]'

# --- Step 3: formatting ---------------------------------------------------
# Column A in the "What's more..." list (rows 46-64) and in the new
# "In Fatigue module" block (rows 66-69) is bold, matching the rest of the
# bullet list. Rows 50-64 already inherited the bold style from the row
# insert above, rows 47 and 66-69 need it applied explicitly.
$ws.Range("A47").Font.Bold = $true
$ws.Range("A66:A69").Font.Bold = $true

# --- Step 4: restore the view -------------------------------------------
$ws.Range("A63").Select()
